# Chapter 6 exercise: build the "Estimates" sheet from the "Data" sheet,
# add Rating/Risk named ranges, and summarise Morningstar ratings / risk.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("Data")

# --- Create "Estimates" as a copy of "Data", placed right after it -------
$data.Copy($null, $data)
$est = $wb.Worksheets.Item(2)
$est.Name = "Estimates"

# --- Named ranges used by the COUNTIF/COUNTA formulas below --------------
$wb.Names.Add('Rating', '=Estimates!$B$2:$B$41')
$wb.Names.Add('Risk', '=Estimates!$C$2:$C$41')

# --- Extra header cells (F1:H1), formatted like the existing header ------
$est.Range("A1").Copy()
$est.Range("F1").PasteSpecial(-4122)
$est.Range("G1").PasteSpecial(-4122)
$est.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$est.Range("F1").Value = "Point Estimate"
$est.Range("G1").Value = "Proportion"

# --- Point-estimate / proportion summary block ----------------------------
$est.Range("E2").Value = "5 Star"
$est.Range("E2").Font.Bold = $true
$est.Range("F2").Formula = '=COUNTIF(Rating, "5 Star")'
$est.Range("G2").Value = 0
$est.Range("G2").NumberFormat = "0%"
$est.Range("G2").Formula = '=F2/F6'

$est.Range("E3").Value = "Above Average"
$est.Range("E3").Font.Bold = $true
$est.Range("F3").Formula = '=COUNTIF(Risk, "Above Average")'
$est.Range("G3").NumberFormat = "0%"
$est.Range("G3").Formula = '=F3/F7'

$est.Range("E4").Value = "<= 2 Star"
$est.Range("E4").Font.Bold = $true
$est.Range("F4").Formula = '=COUNTIF(Rating, "<=2 Star")'
$est.Range("G4").NumberFormat = "0%"
$est.Range("G4").Formula = '=F4/F6'

$est.Range("F5").Value = "Total"
$est.Range("F5").Font.Bold = $true

$est.Range("E6").Value = "Rating"
$est.Range("E6").Font.Bold = $true
$est.Range("F6").Formula = '=COUNTA(Rating)'

$est.Range("E7").Value = "Risk"
$est.Range("E7").Font.Bold = $true
$est.Range("F7").Formula = '=COUNTA(Risk)'

# --- Column widths for the new analysis columns ---------------------------
$est.Columns.Item(4).ColumnWidth = 9.0859375
$est.Columns.Item(5).ColumnWidth = 12.94140625
$est.Columns.Item(6).ColumnWidth = 14.6796875
$est.Columns.Item(7).ColumnWidth = 11.6796875

# --- View state: Data keeps a single-cell selection, Estimates is active --
$data.Range("L1").Select()
$est.Range("H5").Select()
$est.Activate()

$wb.Save()
